# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns
# to the latest scraped snapshot, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '68.120.78'; Numeric = $false },
    @{ Cell = "E2"; Value = '  +0.35%  '; Numeric = $false },
    @{ Cell = "D3"; Value = '3.245.50'; Numeric = $false },
    @{ Cell = "E3"; Value = '  -0.99%  '; Numeric = $false },
    @{ Cell = "E4"; Value = '  +0.13%  '; Numeric = $false },
    @{ Cell = "D5"; Value = '583.41'; Numeric = $true },
    @{ Cell = "E5"; Value = '  +0.35%  '; Numeric = $false },
    @{ Cell = "D6"; Value = '180.68'; Numeric = $true },
    @{ Cell = "E6"; Value = '  -0.99%  '; Numeric = $false },
    @{ Cell = "E7"; Value = '  +0.11%  '; Numeric = $false },
    @{ Cell = "D8"; Value = '0.594'; Numeric = $true },
    @{ Cell = "E8"; Value = '  -1.23%  '; Numeric = $false },
    @{ Cell = "D9"; Value = '0.133'; Numeric = $true },
    @{ Cell = "E9"; Value = '  -0.95%  '; Numeric = $false },
    @{ Cell = "D10"; Value = '6.65'; Numeric = $true },
    @{ Cell = "E10"; Value = '  -1.63%  '; Numeric = $false },
    @{ Cell = "D11"; Value = '0.420'; Numeric = $true },
    @{ Cell = "E11"; Value = '  +0.20%  '; Numeric = $false },
    @{ Cell = "D12"; Value = '3.811.69'; Numeric = $false },
    @{ Cell = "E12"; Value = '  -0.85%  '; Numeric = $false },
    @{ Cell = "D13"; Value = '0.137'; Numeric = $true },
    @{ Cell = "E13"; Value = '  +0.01%  '; Numeric = $false },
    @{ Cell = "D14"; Value = '28.12'; Numeric = $true },
    @{ Cell = "E14"; Value = '  -1.87%  '; Numeric = $false },
    @{ Cell = "D15"; Value = '68.195.23'; Numeric = $false },
    @{ Cell = "E15"; Value = '  +0.56%  '; Numeric = $false },
    @{ Cell = "D16"; Value = '0.0000170'; Numeric = $true },
    @{ Cell = "E16"; Value = '  +0.41%  '; Numeric = $false },
    @{ Cell = "D17"; Value = '3.237.49'; Numeric = $false },
    @{ Cell = "E17"; Value = '  -0.82%  '; Numeric = $false },
    @{ Cell = "D18"; Value = '5.81'; Numeric = $true },
    @{ Cell = "E18"; Value = '  -0.79%  '; Numeric = $false },
    @{ Cell = "D19"; Value = '13.45'; Numeric = $true },
    @{ Cell = "E19"; Value = '  -1.19%  '; Numeric = $false },
    @{ Cell = "D20"; Value = '392.74'; Numeric = $true },
    @{ Cell = "E20"; Value = '  +4.34%  '; Numeric = $false },
    @{ Cell = "D21"; Value = '7.63'; Numeric = $true },
    @{ Cell = "E21"; Value = '  -0.53%  '; Numeric = $false },
    @{ Cell = "D22"; Value = '71.41'; Numeric = $true },
    @{ Cell = "E22"; Value = '  +0.21%  '; Numeric = $false },
    @{ Cell = "D23"; Value = '0.998'; Numeric = $true },
    @{ Cell = "E23"; Value = '  -0.40%  '; Numeric = $false },
    @{ Cell = "D24"; Value = '0.514'; Numeric = $true },
    @{ Cell = "E24"; Value = '  +0.13%  '; Numeric = $false },
    @{ Cell = "D25"; Value = '0.0000118'; Numeric = $true },
    @{ Cell = "E25"; Value = '  -1.85%  '; Numeric = $false },
    @{ Cell = "E26"; Value = '  +4.00%  '; Numeric = $false },
    @{ Cell = "D27"; Value = '9.57'; Numeric = $true },
    @{ Cell = "E27"; Value = '  -0.84%  '; Numeric = $false },
    @{ Cell = "E28"; Value = '  +0.14%  '; Numeric = $false },
    @{ Cell = "E29"; Value = '  -0.22%  '; Numeric = $false },
    @{ Cell = "D30"; Value = '5.66'; Numeric = $true },
    @{ Cell = "E30"; Value = '  -0.96%  '; Numeric = $false },
    @{ Cell = "D31"; Value = '22.95'; Numeric = $true },
    @{ Cell = "E31"; Value = '  +0.57%  '; Numeric = $false },
    @{ Cell = "D32"; Value = '7.10'; Numeric = $true },
    @{ Cell = "E32"; Value = '  +1.91%  '; Numeric = $false },
    @{ Cell = "E34"; Value = '  -1.26%  '; Numeric = $false },
    @{ Cell = "D35"; Value = '164.16'; Numeric = $true },
    @{ Cell = "E35"; Value = '  +0.51%  '; Numeric = $false },
    @{ Cell = "D36"; Value = '1.48'; Numeric = $true },
    @{ Cell = "E36"; Value = '  -1.89%  '; Numeric = $false },
    @{ Cell = "D37"; Value = '1.92'; Numeric = $true },
    @{ Cell = "E37"; Value = '  +3.16%  '; Numeric = $false },
    @{ Cell = "D38"; Value = '0.821'; Numeric = $true },
    @{ Cell = "E38"; Value = '  -3.91%  '; Numeric = $false },
    @{ Cell = "E39"; Value = '  -1.97%  '; Numeric = $false },
    @{ Cell = "D40"; Value = '26.20'; Numeric = $true },
    @{ Cell = "E40"; Value = '  -2.92%  '; Numeric = $false },
    @{ Cell = "D41"; Value = '6.56'; Numeric = $true },
    @{ Cell = "E41"; Value = '  -4.28%  '; Numeric = $false },
    @{ Cell = "D42"; Value = '41.28'; Numeric = $true },
    @{ Cell = "E42"; Value = '  +0.79%  '; Numeric = $false },
    @{ Cell = "D43"; Value = '0.0688'; Numeric = $true },
    @{ Cell = "E43"; Value = '  +0.44%  '; Numeric = $false },
    @{ Cell = "D44"; Value = '2.46'; Numeric = $true },
    @{ Cell = "E44"; Value = '  -6.45%  '; Numeric = $false },
    @{ Cell = "D45"; Value = '340.96'; Numeric = $true },
    @{ Cell = "E45"; Value = '  -4.26%  '; Numeric = $false },
    @{ Cell = "D46"; Value = '2.580.72'; Numeric = $false },
    @{ Cell = "E46"; Value = '  -4.88%  '; Numeric = $false },
    @{ Cell = "D47"; Value = '24.58'; Numeric = $true },
    @{ Cell = "E47"; Value = '  -3.72%  '; Numeric = $false },
    @{ Cell = "E48"; Value = '  -0.53%  '; Numeric = $false },
    @{ Cell = "D49"; Value = '31.68'; Numeric = $true },
    @{ Cell = "E49"; Value = '  +0.68%  '; Numeric = $false },
    @{ Cell = "E50"; Value = '  +1.73%  '; Numeric = $false },
    @{ Cell = "D51"; Value = '0.101'; Numeric = $true },
    @{ Cell = "E51"; Value = '  -1.37%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Value looks like a plain number ("583.41"); prefix with an
        # apostrophe so Excel stores it as literal text (matches the
        # source workbook, where Price is a text column), then drop the
        # auto-added "number stored as text" formatting it triggers so the
        # cell keeps the sheet's default (unstyled) look.
        $range.Value = "'" + $u.Value
        $range.ClearFormats()
    } else {
        $range.Value = $u.Value
    }
}
